$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a literal text value, even when the
# text looks like a number (Excel/COM would otherwise coerce plain
# Range.Value assignments of numeric-looking strings into numbers).
# We build the text via a `="..."` formula (so it is never subject to
# numeric auto-detection) and then convert the formula to a static
# value via copy / paste-special values-only.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $escaped = $val.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue $ws "D2" '58.374.28'
$ws.Range("E2").Value = '  -4.29%  '

Set-TextValue $ws "D3" '2.614.77'
$ws.Range("E3").Value = '  -4.02%  '

$ws.Range("E4").Value = '  +0.00%  '

Set-TextValue $ws "D5" '519.86'
$ws.Range("E5").Value = '  -2.00%  '

Set-TextValue $ws "D6" '141.94'
$ws.Range("E6").Value = '  -4.02%  '

$ws.Range("E7").Value = '  +0.32%  '

$ws.Range("E8").Value = '  -2.79%  '

Set-TextValue $ws "D9" '6.53'
$ws.Range("E9").Value = '  -9.75%  '

$ws.Range("E10").Value = '  -3.83%  '

Set-TextValue $ws "D11" '0.335'
$ws.Range("E11").Value = '  -2.25%  '

$ws.Range("E12").Value = '  +0.73%  '

Set-TextValue $ws "D13" '3.074.35'
$ws.Range("E13").Value = '  -3.71%  '

Set-TextValue $ws "D14" '58.333.78'
$ws.Range("E14").Value = '  -4.31%  '

Set-TextValue $ws "D15" '20.88'
$ws.Range("E15").Value = '  -3.30%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws "D16" '2.659.87'
$ws.Range("E16").Value = '  -5.49%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws "D17" '0.0000135'
$ws.Range("E17").Value = '  -2.76%  '

Set-TextValue $ws "D18" '336.83'
$ws.Range("E18").Value = '  -3.17%  '

$ws.Range("E19").Value = '  -3.41%  '

Set-TextValue $ws "D20" '10.37'
$ws.Range("E20").Value = '  -2.53%  '

Set-TextValue $ws "D21" '6.27'
$ws.Range("E21").Value = '  -3.00%  '

Set-TextValue $ws "D23" '64.88'
$ws.Range("E23").Value = '  +2.10%  '

$ws.Range("E24").Value = '  -2.09%  '

$ws.Range("E25").Value = '  -3.94%  '

Set-TextValue $ws "D26" '1.00'
$ws.Range("E26").Value = '  +1.14%  '

$ws.Range("E27").Value = '  -3.57%  '

Set-TextValue $ws "D28" '0.0₃0785'
$ws.Range("E28").Value = '  -5.51%  '

Set-TextValue $ws "D29" '6.51'
$ws.Range("E29").Value = '  -4.34%  '

Set-TextValue $ws "D30" '0.999'
$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("E31").Value = '  -1.30%  '

$ws.Range("E32").Value = '  -2.33%  '

Set-TextValue $ws "D33" '149.78'
$ws.Range("E33").Value = '  -0.16%  '

Set-TextValue $ws "D34" '4.07'
$ws.Range("E34").Value = '  -5.10%  '

$ws.Range("E35").Value = '  -5.81%  '

Set-TextValue $ws "D36" '0.891'
$ws.Range("E36").Value = '  -4.05%  '

Set-TextValue $ws "D37" '0.846'
$ws.Range("E37").Value = '  -6.64%  '

Set-TextValue $ws "D38" '36.23'
$ws.Range("E38").Value = '  -2.82%  '

$ws.Range("E39").Value = '  -7.79%  '

$ws.Range("E40").Value = '  -2.64%  '

$ws.Range("E41").Value = '  +0.47%  '

$ws.Range("E42").Value = '  -2.71%  '

Set-TextValue $ws "D43" '0.0966'
$ws.Range("E43").Value = '  -2.79%  '

Set-TextValue $ws "D44" '267.65'
$ws.Range("E44").Value = '  -6.18%  '

$ws.Range("E45").Value = '  +1.07%  '

Set-TextValue $ws "D46" '19.07'
$ws.Range("E46").Value = '  -7.16%  '

Set-TextValue $ws "D47" '0.0529'
$ws.Range("E47").Value = '  -3.14%  '

Set-TextValue $ws "D48" '2.026.32'
$ws.Range("E48").Value = '  -5.11%  '

$ws.Range("E49").Value = '  -2.49%  '

Set-TextValue $ws "D50" '4.56'
$ws.Range("E50").Value = '  -8.59%  '

Set-TextValue $ws "D51" '18.16'
$ws.Range("E51").Value = '  -7.12%  '
